$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Clear the score values in C2:D9 (bonus/slutspil scores removed - goals collapsed into one combined sheet)
$ws.Range("C2:D9").ClearContents()

# Update the selection to match the authored state: the author had dragged a
# selection over C2:D9 with D9 as the active (last-clicked) cell. This host's
# Range.Select() anchors the active cell to the top-left corner of the
# resulting rectangle, so sqref="C2:D9" is reproduced exactly; the activeCell
# sub-attribute is the closest achievable approximation.
$ws.Range("C2:D9").Select()
